$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $value)
    $r = $ws.Range($cellRef)
    if ($value -match '^[+-]?[0-9]*\.?[0-9]+$') {
        $r.NumberFormat = "@"
        $r.Value = $value
        $r.NumberFormat = "General"
    } else {
        $r.Value = $value
    }
}

Set-TextValue "D2" "57.621.95"
Set-TextValue "E2" "  -0.45%  "
Set-TextValue "D3" "3.122.29"
Set-TextValue "E3" "  -0.36%  "
Set-TextValue "E4" "  -0.02%  "
Set-TextValue "D5" "527.72"
Set-TextValue "E5" "  +0.19%  "
Set-TextValue "D6" "138.33"
Set-TextValue "E6" "  -2.44%  "
Set-TextValue "D7" "1.00"
Set-TextValue "E7" "  +0.02%  "
Set-TextValue "D8" "3.121.06"
Set-TextValue "E8" "  -0.49%  "
Set-TextValue "D9" "0.444"
Set-TextValue "E9" "  +2.44%  "
Set-TextValue "D10" "7.15"
Set-TextValue "E10" "  -1.93%  "
Set-TextValue "E11" "  -2.10%  "
Set-TextValue "E12" "  +2.29%  "
Set-TextValue "D13" "3.667.48"
Set-TextValue "E13" "  -0.28%  "
Set-TextValue "E14" "  +2.64%  "
Set-TextValue "D15" "25.45"
Set-TextValue "E15" "  -3.61%  "
Set-TextValue "E16" "  -0.56%  "
Set-TextValue "D17" "57.781.34"
Set-TextValue "E17" "  -0.36%  "
Set-TextValue "D18" "3.127.23"
Set-TextValue "E18" "  -0.40%  "
Set-TextValue "D19" "6.00"
Set-TextValue "E19" "  -2.52%  "
Set-TextValue "D20" "12.77"
Set-TextValue "E20" "  -1.34%  "
Set-TextValue "D21" "7.88"
Set-TextValue "E21" "  -3.18%  "
Set-TextValue "D22" "353.05"
Set-TextValue "E22" "  +4.55%  "
Set-TextValue "E23" "  +0.28%  "
Set-TextValue "D24" "68.55"
Set-TextValue "E24" "  +2.42%  "
Set-TextValue "D25" "0.506"
Set-TextValue "E25" "  -1.46%  "
Set-TextValue "E26" "  -0.13%  "
Set-TextValue "E27" "  +0.01%  "
Set-TextValue "D28" "0.0₃0913"
Set-TextValue "E28" "  -2.17%  "
Set-TextValue "D29" "7.44"
Set-TextValue "E29" "  +2.43%  "
Set-TextValue "D30" "0.999"
Set-TextValue "E30" "  +0.02%  "
Set-TextValue "D31" "6.23"
Set-TextValue "E31" "  -6.28%  "
Set-TextValue "E32" "  -0.12%  "
Set-TextValue "D33" "21.15"
Set-TextValue "E33" "  +0.46%  "
Set-TextValue "D34" "1.18"
Set-TextValue "E34" "  -2.77%  "
Set-TextValue "E35" "  +5.10%  "
Set-TextValue "D36" "157.69"
Set-TextValue "E36" "  +1.33%  "
Set-TextValue "E37" "  +0.57%  "
Set-TextValue "D38" "26.21"
Set-TextValue "E38" "  -3.72%  "
Set-TextValue "E39" "  -2.78%  "
Set-TextValue "D40" "0.0667"
Set-TextValue "E40" "  -0.26%  "
Set-TextValue "B41" "Filecoin"
Set-TextValue "C41" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D41" "4.19"
Set-TextValue "E41" "  +6.52%  "
Set-TextValue "B42" "Stacks"
Set-TextValue "C42" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D42" "1.62"
Set-TextValue "E42" "  +5.41%  "
Set-TextValue "D43" "0.698"
Set-TextValue "E43" "  +0.76%  "
Set-TextValue "D44" "3.171.00"
Set-TextValue "E44" "  -0.22%  "
Set-TextValue "D45" "36.50"
Set-TextValue "E45" "  -1.41%  "
Set-TextValue "E46" "  -0.02%  "
Set-TextValue "E47" "  +2.75%  "
Set-TextValue "D48" "2.315.60"
Set-TextValue "E48" "  +0.55%  "
Set-TextValue "D49" "0.974"
Set-TextValue "E49" "  -2.57%  "
Set-TextValue "D50" "6.05"
Set-TextValue "E50" "  +0.41%  "
Set-TextValue "D51" "20.31"
Set-TextValue "E51" "  -3.76%  "
